$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados..." timestamp in A1
$ws.Range("A1").Value = 'Datos actualizados a 23 de Marzo de 2020 a las 09:46'

# Row-level updates: country name (col A) + Casos totales/Nuevos casos/Casos
# activas/Recuperados/Casos criticos/Muertes hoy/Muertes (cols B-H).
# Rows shift their country label because the underlying shared-string table
# was re-ordered (ranking changed) in addition to fresh case counts.
$rowData = @(
    @{ Row = 6; Values = @('Estados Unidos', 35070, 1524, 178, 34434, 795, 39, 458) },
    @{ Row = 8; Values = @('Alemania', 24904, 31, 266, 24544, 23, 0, 94) },
    @{ Row = 20; Values = @('Brasil', 1604, 58, 2, 1577, 18, 0, 25) },
    @{ Row = 21; Values = @('Portugal', 1600, 0, 5, 1581, 26, 0, 14) },
    @{ Row = 43; Values = @('Filipinas', 462, 82, 18, 411, 1, 8, 33) },
    @{ Row = 44; Values = @('Singapur', 455, 0, 144, 309, 14, 0, 2) },
    @{ Row = 45; Values = @('Rusia', 438, 71, 16, 421, 0, 0, 1) },
    @{ Row = 46; Values = @('Rumania', 433, 0, 64, 366, 14, 0, 3) },
    @{ Row = 47; Values = @('India', 425, 29, 24, 393, 0, 1, 8) },
    @{ Row = 48; Values = @('Eslovenia', 414, 0, 0, 412, 12, 0, 2) },
    @{ Row = 50; Values = @('Barein', 337, 3, 160, 175, 3, 0, 2) },
    @{ Row = 56; Values = @('Croacia', 306, 52, 5, 300, 5, 0, 1) },
    @{ Row = 57; Values = @('Sudafrica', 274, 0, 2, 272, 0, 0, 0) },
    @{ Row = 58; Values = @('Argentina', 266, 0, 27, 235, 0, 0, 4) },
    @{ Row = 90; Values = @('Sri Lanka', 87, 5, 3, 84, 2, 0, 0) },
    @{ Row = 91; Values = @('Camboya', 86, 2, 2, 84, 1, 0, 0) },
    @{ Row = 168; Values = @('Niger', 2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 169; Values = @('Benin', 2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 170; Values = @('Butan', 2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 171; Values = @('Guinea', 2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 173; Values = @('Angola', 2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 175; Values = @('Haiti', 2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 176; Values = @('Nicaragua', 2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 177; Values = @('Sudan', 2, 0, 0, 1, 0, 0, 1) },
    @{ Row = 178; Values = @('Nepal', 2, 1, 1, 1, 0, 0, 0) },
    @{ Row = 179; Values = @('Somalia', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 180; Values = @('Republica de Yibuti', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 181; Values = @('Siria', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 182; Values = @('San Martin (Parte Holandesa)', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 183; Values = @('Montserrat', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 184; Values = @('Antigua y Barbuda', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 185; Values = @('Republica del Chad', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 186; Values = @('Dominica', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 187; Values = @('Granada', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 188; Values = @('Santa Sede', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 189; Values = @('San Vicente y las Granadinas', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 190; Values = @('Papua Nueva Guinea', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 191; Values = @('Timor Oriental', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 192; Values = @('Eritrea', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 193; Values = @('Uganda', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 194; Values = @('Mozambique', 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 195; Values = @('Gambia', 1, 0, 0, 0, 0, 1, 1) }
)

foreach ($entry in $rowData) {
    $rowNum = $entry.Row
    $vals = $entry.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($rowNum, $i + 1).Value = $vals[$i]
    }
}